$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Range("E2").Value = 3
    $ws.Range("G2").Value = 110.8604276666667
    $ws.Range("H2").Value = 332.581283
    $ws.Range("I2").Value = 0.2509786052589675
    $ws.Range("J2").Value = 0.2509786052589675
    $ws.Range("K2").Value = 3
    $ws.Range("M2").Value = 22.839587
    $ws.Range("N2").Value = 68.518761
    $ws.Range("O2").Value = 0.2024156068965367
    $ws.Range("P2").Value = 0.2024156068965367
    $ws.Range("Q2").Value = 2532.00638255004
    $ws.Range("R2").Value = 22788.05744295036
    $ws.Range("S2").Value = 0.05080198670154022
    $ws.Range("T2").Value = 0.05080198670154022
    $ws.Range("E3").Value = 3
    $ws.Range("G3").Value = 110.8604276666667
    $ws.Range("H3").Value = 332.581283
    $ws.Range("I3").Value = 0.2509786052589675
    $ws.Range("J3").Value = 0.2509786052589675
    $ws.Range("K3").Value = 3
    $ws.Range("M3").Value = 31.29092
    $ws.Range("N3").Value = 93.87276
    $ws.Range("O3").Value = 0.2773154594323872
    $ws.Range("P3").Value = 0.2773154594323872
    $ws.Range("Q3").Value = 3468.924773283453
    $ws.Range("R3").Value = 31220.32295955108
    $ws.Range("S3").Value = 0.06960024722509031
    $ws.Range("T3").Value = 0.06960024722509033
    $ws.Range("E4").Value = 3
    $ws.Range("G4").Value = 110.8604276666667
    $ws.Range("H4").Value = 332.581283
    $ws.Range("I4").Value = 0.2509786052589675
    $ws.Range("J4").Value = 0.2509786052589675
    $ws.Range("K4").Value = 3
    $ws.Range("M4").Value = 35.824351
    $ws.Range("N4").Value = 107.473053
    $ws.Range("O4").Value = 0.3174929454433458
    $ws.Range("P4").Value = 0.3174929454433459
    $ws.Range("Q4").Value = 3971.502872740778
    $ws.Range("R4").Value = 35743.525854667
    $ws.Range("S4").Value = 0.07968393662693239
    $ws.Range("T4").Value = 0.07968393662693241
    $ws.Range("E5").Value = 3
    $ws.Range("G5").Value = 110.8604276666667
    $ws.Range("H5").Value = 332.581283
    $ws.Range("I5").Value = 0.2509786052589675
    $ws.Range("J5").Value = 0.2509786052589675
    $ws.Range("K5").Value = 3
    $ws.Range("M5").Value = 22.88025066666667
    $ws.Range("N5").Value = 68.64075199999999
    $ws.Range("O5").Value = 0.2027759882277303
    $ws.Range("P5").Value = 0.2027759882277303
    $ws.Range("Q5").Value = 2536.514374027201
    $ws.Range("R5").Value = 22828.62936624481
    $ws.Range("S5").Value = 0.05089243470540455
    $ws.Range("T5").Value = 0.05089243470540455
    $ws.Range("E6").Value = 3
    $ws.Range("G6").Value = 184.841802
    $ws.Range("H6").Value = 554.525406
    $ws.Range("I6").Value = 0.4184661617850055
    $ws.Range("J6").Value = 0.4184661617850055
    $ws.Range("K6").Value = 3
    $ws.Range("M6").Value = 22.839587
    $ws.Range("N6").Value = 68.518761
    $ws.Range("O6").Value = 0.2024156068965367
    $ws.Range("P6").Value = 0.2024156068965367
    $ws.Range("Q6").Value = 4221.710418015774
    $ws.Range("R6").Value = 37995.39376214196
    $ws.Range("S6").Value = 0.08470408210337618
    $ws.Range("T6").Value = 0.0847040821033762
    $ws.Range("E7").Value = 3
    $ws.Range("G7").Value = 184.841802
    $ws.Range("H7").Value = 554.525406
    $ws.Range("I7").Value = 0.4184661617850055
    $ws.Range("J7").Value = 0.4184661617850055
    $ws.Range("K7").Value = 3
    $ws.Range("M7").Value = 31.29092
    $ws.Range("N7").Value = 93.87276
    $ws.Range("O7").Value = 0.2773154594323872
    $ws.Range("P7").Value = 0.2773154594323872
    $ws.Range("Q7").Value = 5783.87003903784
    $ws.Range("R7").Value = 52054.83035134056
    $ws.Range("S7").Value = 0.1160471359123164
    $ws.Range("T7").Value = 0.1160471359123165
    $ws.Range("E8").Value = 3
    $ws.Range("G8").Value = 184.841802
    $ws.Range("H8").Value = 554.525406
    $ws.Range("I8").Value = 0.4184661617850055
    $ws.Range("J8").Value = 0.4184661617850055
    $ws.Range("K8").Value = 3
    $ws.Range("M8").Value = 35.824351
    $ws.Range("N8").Value = 107.473053
    $ws.Range("O8").Value = 0.3174929454433458
    $ws.Range("P8").Value = 0.3174929454433459
    $ws.Range("Q8").Value = 6621.837594320502
    $ws.Range("R8").Value = 59596.53834888452
    $ws.Range("S8").Value = 0.1328600542734931
    $ws.Range("T8").Value = 0.1328600542734931
    $ws.Range("E9").Value = 3
    $ws.Range("G9").Value = 184.841802
    $ws.Range("H9").Value = 554.525406
    $ws.Range("I9").Value = 0.4184661617850055
    $ws.Range("J9").Value = 0.4184661617850055
    $ws.Range("K9").Value = 3
    $ws.Range("M9").Value = 22.88025066666667
    $ws.Range("N9").Value = 68.64075199999999
    $ws.Range("O9").Value = 0.2027759882277303
    $ws.Range("P9").Value = 0.2027759882277303
    $ws.Range("Q9").Value = 4229.226763438368
    $ws.Range("R9").Value = 38063.0408709453
    $ws.Range("S9").Value = 0.08485488949581974
    $ws.Range("T9").Value = 0.08485488949581974
    $ws.Range("E10").Value = 3
    $ws.Range("G10").Value = 95.23175666666667
    $ws.Range("H10").Value = 285.69527
    $ws.Range("I10").Value = 0.2155966197102082
    $ws.Range("J10").Value = 0.2155966197102082
    $ws.Range("K10").Value = 3
    $ws.Range("M10").Value = 22.839587
    $ws.Range("N10").Value = 68.518761
    $ws.Range("O10").Value = 0.2024156068965367
    $ws.Range("P10").Value = 0.2024156068965367
    $ws.Range("Q10").Value = 2175.053991551163
    $ws.Range("R10").Value = 19575.48592396047
    $ws.Range("S10").Value = 0.04364012062348362
    $ws.Range("T10").Value = 0.04364012062348363
    $ws.Range("E11").Value = 3
    $ws.Range("G11").Value = 95.23175666666667
    $ws.Range("H11").Value = 285.69527
    $ws.Range("I11").Value = 0.2155966197102082
    $ws.Range("J11").Value = 0.2155966197102082
    $ws.Range("K11").Value = 3
    $ws.Range("M11").Value = 31.29092
    $ws.Range("N11").Value = 93.87276
    $ws.Range("O11").Value = 0.2773154594323872
    $ws.Range("P11").Value = 0.2773154594323872
    $ws.Range("Q11").Value = 2979.889279316134
    $ws.Range("R11").Value = 26819.0035138452
    $ws.Range("S11").Value = 0.05978827564700605
    $ws.Range("T11").Value = 0.05978827564700607
    $ws.Range("E12").Value = 3
    $ws.Range("G12").Value = 95.23175666666667
    $ws.Range("H12").Value = 285.69527
    $ws.Range("I12").Value = 0.2155966197102082
    $ws.Range("J12").Value = 0.2155966197102082
    $ws.Range("K12").Value = 3
    $ws.Range("M12").Value = 35.824351
    $ws.Range("N12").Value = 107.473053
    $ws.Range("O12").Value = 0.3174929454433458
    $ws.Range("P12").Value = 0.3174929454433459
    $ws.Range("Q12").Value = 3411.615877173257
    $ws.Range("R12").Value = 30704.54289455931
    $ws.Range("S12").Value = 0.06845040581942291
    $ws.Range("T12").Value = 0.06845040581942294
    $ws.Range("E13").Value = 3
    $ws.Range("G13").Value = 95.23175666666667
    $ws.Range("H13").Value = 285.69527
    $ws.Range("I13").Value = 0.2155966197102082
    $ws.Range("J13").Value = 0.2155966197102082
    $ws.Range("K13").Value = 3
    $ws.Range("M13").Value = 22.88025066666667
    $ws.Range("N13").Value = 68.64075199999999
    $ws.Range("O13").Value = 0.2027759882277303
    $ws.Range("P13").Value = 0.2027759882277303
    $ws.Range("Q13").Value = 2178.926463960338
    $ws.Range("R13").Value = 19610.33817564304
    $ws.Range("S13").Value = 0.04371781762029562
    $ws.Range("T13").Value = 0.04371781762029563
    $ws.Range("E14").Value = 3
    $ws.Range("G14").Value = 50.778675
    $ws.Range("H14").Value = 152.336025
    $ws.Range("I14").Value = 0.1149586132458188
    $ws.Range("J14").Value = 0.1149586132458188
    $ws.Range("K14").Value = 3
    $ws.Range("M14").Value = 22.839587
    $ws.Range("N14").Value = 68.518761
    $ws.Range("O14").Value = 0.2024156068965367
    $ws.Range("P14").Value = 0.2024156068965367
    $ws.Range("Q14").Value = 1159.763965407225
    $ws.Range("R14").Value = 10437.87568866502
    $ws.Range("S14").Value = 0.02326941746813665
    $ws.Range("T14").Value = 0.02326941746813666
    $ws.Range("E15").Value = 3
    $ws.Range("G15").Value = 50.778675
    $ws.Range("H15").Value = 152.336025
    $ws.Range("I15").Value = 0.1149586132458188
    $ws.Range("J15").Value = 0.1149586132458188
    $ws.Range("K15").Value = 3
    $ws.Range("M15").Value = 31.29092
    $ws.Range("N15").Value = 93.87276
    $ws.Range("O15").Value = 0.2773154594323872
    $ws.Range("P15").Value = 0.2773154594323872
    $ws.Range("Q15").Value = 1588.911457131
    $ws.Range("R15").Value = 14300.203114179
    $ws.Range("S15").Value = 0.03187980064797435
    $ws.Range("T15").Value = 0.03187980064797435
    $ws.Range("E16").Value = 3
    $ws.Range("G16").Value = 50.778675
    $ws.Range("H16").Value = 152.336025
    $ws.Range("I16").Value = 0.1149586132458188
    $ws.Range("J16").Value = 0.1149586132458188
    $ws.Range("K16").Value = 3
    $ws.Range("M16").Value = 35.824351
    $ws.Range("N16").Value = 107.473053
    $ws.Range("O16").Value = 0.3174929454433458
    $ws.Range("P16").Value = 0.3174929454433459
    $ws.Range("Q16").Value = 1819.113076514925
    $ws.Range("R16").Value = 16372.01768863433
    $ws.Range("S16").Value = 0.03649854872349743
    $ws.Range("T16").Value = 0.03649854872349745
    $ws.Range("E17").Value = 3
    $ws.Range("G17").Value = 50.778675
    $ws.Range("H17").Value = 152.336025
    $ws.Range("I17").Value = 0.1149586132458188
    $ws.Range("J17").Value = 0.1149586132458188
    $ws.Range("K17").Value = 3
    $ws.Range("M17").Value = 22.88025066666667
    $ws.Range("N17").Value = 68.64075199999999
    $ws.Range("O17").Value = 0.2027759882277303
    $ws.Range("P17").Value = 0.2027759882277303
    $ws.Range("Q17").Value = 1161.8288125212
    $ws.Range("R17").Value = 10456.4593126908
    $ws.Range("S17").Value = 0.02331084640621035
    $ws.Range("T17").Value = 0.02331084640621035

Write-Output "Updated $($ws.Name) with revised expressing-cell counts and derived metrics."
